$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Cells.Item(16, 1).Value = 111814104
$ws.Cells.Item(16, 2).Value = 56398
$ws.Cells.Item(16, 5).Value = 100109
$ws.Cells.Item(16, 6).Value = 'Tretåig hackspett'
$ws.Cells.Item(16, 7).Value = 'Picoides tridactylus'
$ws.Cells.Item(16, 8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(16, 17).Value = 610154.5078508666
$ws.Cells.Item(16, 18).Value = 7121460.305022033
$ws.Cells.Item(16, 26).Value = '17:23'
$ws.Cells.Item(16, 28).Value = '17:23'

# Row 17
$ws.Cells.Item(17, 1).Value = 111815269
$ws.Cells.Item(17, 2).Value = 90666
$ws.Cells.Item(17, 4).Value = 'LC'
$ws.Cells.Item(17, 5).Value = 4364
$ws.Cells.Item(17, 6).Value = 'Dropptaggsvamp'
$ws.Cells.Item(17, 7).Value = 'Hydnellum ferrugineum'
$ws.Cells.Item(17, 8).Value = '(Fr.:Fr.) P. Karst.'
$ws.Cells.Item(17, 17).Value = 610053.7842541422
$ws.Cells.Item(17, 18).Value = 7121273.15248157
$ws.Cells.Item(17, 26).Value = '18:27'
$ws.Cells.Item(17, 28).Value = '18:27'

# Row 18
$ws.Cells.Item(18, 1).Value = 111815024
$ws.Cells.Item(18, 2).Value = 56414
$ws.Cells.Item(18, 4).Value = 'NT'
$ws.Cells.Item(18, 5).Value = 100049
$ws.Cells.Item(18, 6).Value = 'Spillkråka'
$ws.Cells.Item(18, 7).Value = 'Dryocopus martius'
$ws.Cells.Item(18, 8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(18, 17).Value = 609922.1399673244
$ws.Cells.Item(18, 18).Value = 7121488.212810148
$ws.Cells.Item(18, 26).Value = '18:12'
$ws.Cells.Item(18, 28).Value = '18:12'

# Row 19
$ws.Cells.Item(19, 1).Value = 111815114
$ws.Cells.Item(19, 2).Value = 90660
$ws.Cells.Item(19, 4).Value = 'NT'
$ws.Cells.Item(19, 5).Value = 4362
$ws.Cells.Item(19, 6).Value = 'Blå taggsvamp'
$ws.Cells.Item(19, 7).Value = 'Hydnellum caeruleum'
$ws.Cells.Item(19, 8).Value = '(Hornem.) P.Karst.'
$ws.Cells.Item(19, 16).Value = 'åsele 1:1, Ås lm'
$ws.Cells.Item(19, 17).Value = 610384.0265214761
$ws.Cells.Item(19, 18).Value = 7121170.261031131
$ws.Cells.Item(19, 19).Value = 5
$ws.Cells.Item(19, 26).Value = '18:19'
$ws.Cells.Item(19, 28).Value = '18:19'

# Row 20
$ws.Cells.Item(20, 1).Value = 111814688
$ws.Cells.Item(20, 2).Value = 90087
$ws.Cells.Item(20, 4).Value = 'LC'
$ws.Cells.Item(20, 5).Value = 3298
$ws.Cells.Item(20, 6).Value = 'Trådticka'
$ws.Cells.Item(20, 7).Value = 'Climacocystis borealis'
$ws.Cells.Item(20, 8).Value = '(Fr.) Kotl. & Pouzar'
$ws.Cells.Item(20, 17).Value = 610011.2059644217
$ws.Cells.Item(20, 18).Value = 7121475.688616944
$ws.Cells.Item(20, 26).Value = '17:55'
$ws.Cells.Item(20, 28).Value = '17:55'

# Row 21
$ws.Cells.Item(21, 1).Value = 111814591
$ws.Cells.Item(21, 2).Value = 77515
$ws.Cells.Item(21, 5).Value = 6425
$ws.Cells.Item(21, 6).Value = 'Garnlav'
$ws.Cells.Item(21, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(21, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(21, 17).Value = 610012.4812897337
$ws.Cells.Item(21, 18).Value = 7121464.398116477
$ws.Cells.Item(21, 26).Value = '17:50'
$ws.Cells.Item(21, 28).Value = '17:50'

# Row 22
$ws.Cells.Item(22, 1).Value = 111814478
$ws.Cells.Item(22, 2).Value = 77515
$ws.Cells.Item(22, 5).Value = 6425
$ws.Cells.Item(22, 6).Value = 'Garnlav'
$ws.Cells.Item(22, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(22, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(22, 16).Value = 'åsele 1:1 (åsele 1:1), Ås lm'
$ws.Cells.Item(22, 17).Value = 610155.3487898401
$ws.Cells.Item(22, 18).Value = 7121461.207019502
$ws.Cells.Item(22, 19).Value = 1
$ws.Cells.Item(22, 26).Value = '17:41'
$ws.Cells.Item(22, 28).Value = '17:41'

# Row 25
$ws.Cells.Item(25, 1).Value = 112013697
$ws.Cells.Item(25, 2).Value = 89423
$ws.Cells.Item(25, 5).Value = 5432
$ws.Cells.Item(25, 6).Value = 'Granticka'
$ws.Cells.Item(25, 7).Value = 'Porodaedalea chrysoloma'
$ws.Cells.Item(25, 8).Value = '(Fr.) Fiasson & Niemelä'
$ws.Cells.Item(25, 17).Value = 610102.0736959254
$ws.Cells.Item(25, 18).Value = 7121412.654772604
$ws.Cells.Item(25, 26).Value = '19:35'
$ws.Cells.Item(25, 28).Value = '19:35'

# Row 26
$ws.Cells.Item(26, 1).Value = 112013691
$ws.Cells.Item(26, 2).Value = 88489
$ws.Cells.Item(26, 5).Value = 1962
$ws.Cells.Item(26, 6).Value = 'Vaddporing'
$ws.Cells.Item(26, 7).Value = 'Anomoporia kamtschatica'
$ws.Cells.Item(26, 8).Value = '(Parmasto) Bondartseva'
$ws.Cells.Item(26, 17).Value = 610134.4051595986
$ws.Cells.Item(26, 18).Value = 7121460.896015909
$ws.Cells.Item(26, 26).Value = '19:29'
$ws.Cells.Item(26, 28).Value = '19:29'

# Row 27
$ws.Cells.Item(27, 1).Value = 112013700
$ws.Cells.Item(27, 2).Value = 77515
$ws.Cells.Item(27, 5).Value = 6425
$ws.Cells.Item(27, 6).Value = 'Garnlav'
$ws.Cells.Item(27, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(27, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(27, 17).Value = 610101.9650201321
$ws.Cells.Item(27, 18).Value = 7121415.702941997

# Row 28
$ws.Cells.Item(28, 1).Value = 112013696
$ws.Cells.Item(28, 2).Value = 86961
$ws.Cells.Item(28, 5).Value = 4962
$ws.Cells.Item(28, 6).Value = 'Mjölsvärting'
$ws.Cells.Item(28, 7).Value = 'Lyophyllum semitale'
$ws.Cells.Item(28, 8).Value = '(Fr. : Fr.) Kühner'
$ws.Cells.Item(28, 17).Value = 610070.1349689787
$ws.Cells.Item(28, 18).Value = 7121402.360087069
$ws.Cells.Item(28, 26).Value = '19:40'
$ws.Cells.Item(28, 28).Value = '19:40'

# Row 29
$ws.Cells.Item(29, 1).Value = 112013703
$ws.Cells.Item(29, 2).Value = 77515
$ws.Cells.Item(29, 5).Value = 6425
$ws.Cells.Item(29, 6).Value = 'Garnlav'
$ws.Cells.Item(29, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(29, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(29, 17).Value = 610144.4332068264
$ws.Cells.Item(29, 18).Value = 7121461.253672058
$ws.Cells.Item(29, 26).Value = '19:28'
$ws.Cells.Item(29, 28).Value = '19:28'

# Row 30
$ws.Cells.Item(30, 1).Value = 112013698
$ws.Cells.Item(30, 2).Value = 77515
$ws.Cells.Item(30, 5).Value = 6425
$ws.Cells.Item(30, 6).Value = 'Garnlav'
$ws.Cells.Item(30, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(30, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(30, 17).Value = 610094.4326785516
$ws.Cells.Item(30, 18).Value = 7121455.546697079
$ws.Cells.Item(30, 26).Value = '19:49'
$ws.Cells.Item(30, 28).Value = '19:49'

# Row 31
$ws.Cells.Item(31, 1).Value = 112013690
$ws.Cells.Item(31, 2).Value = 88489
$ws.Cells.Item(31, 5).Value = 1962
$ws.Cells.Item(31, 6).Value = 'Vaddporing'
$ws.Cells.Item(31, 7).Value = 'Anomoporia kamtschatica'
$ws.Cells.Item(31, 8).Value = '(Parmasto) Bondartseva'
$ws.Cells.Item(31, 17).Value = 610051.8565798617
$ws.Cells.Item(31, 18).Value = 7121425.252971379
$ws.Cells.Item(31, 26).Value = '19:43'
$ws.Cells.Item(31, 28).Value = '19:43'

# Row 32
$ws.Cells.Item(32, 1).Value = 112013704
$ws.Cells.Item(32, 2).Value = 81248
$ws.Cells.Item(32, 5).Value = 1312
$ws.Cells.Item(32, 6).Value = 'Gammelgransskål'
$ws.Cells.Item(32, 7).Value = 'Pseudographis pinicola'
$ws.Cells.Item(32, 8).Value = '(Nyl.) Rehm'
$ws.Cells.Item(32, 17).Value = 610093.591720929
$ws.Cells.Item(32, 18).Value = 7121454.644715369
$ws.Cells.Item(32, 26).Value = '19:49'
$ws.Cells.Item(32, 28).Value = '19:49'

# Row 33
$ws.Cells.Item(33, 1).Value = 112013699
$ws.Cells.Item(33, 17).Value = 610068.1736430819
$ws.Cells.Item(33, 18).Value = 7121408.394281525
$ws.Cells.Item(33, 26).Value = '19:40'
$ws.Cells.Item(33, 28).Value = '19:40'
